$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay plain text (as in the source data) so
# numeric-looking strings like "1.00" or "20.90" are not silently coerced
# into numbers (which would drop the trailing zero) when assigned below.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '37.767.66'
$ws.Range("E2").Value = '  -0.07%  '

# Row 3
$ws.Range("D3").Value = '2.079.64'
$ws.Range("E3").Value = '  -1.14%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").Value = '233.78'
$ws.Range("E5").Value = '  -0.34%  '

# Row 6
$ws.Range("D6").Value = '0.624'
$ws.Range("E6").Value = '  +0.22%  '

# Row 7
$ws.Range("B7").Value = 'Solana'
$ws.Range("C7").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D7").Value = '58.55'
$ws.Range("E7").Value = '  +1.06%  '

# Row 8
$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.01%  '

# Row 9
$ws.Range("D9").Value = '0.393'
$ws.Range("E9").Value = '  +0.74%  '

# Row 10
$ws.Range("E10").Value = '  +0.89%  '

# Row 11
$ws.Range("E11").Value = '  +2.95%  '

# Row 13
$ws.Range("D13").Value = '14.79'
$ws.Range("E13").Value = '  +2.06%  '

# Row 14
$ws.Range("D14").Value = '20.90'
$ws.Range("E14").Value = '  -2.59%  '

# Row 15
$ws.Range("E15").Value = '  -0.80%  '

# Row 16
$ws.Range("D16").Value = '5.33'
$ws.Range("E16").Value = '  +2.29%  '

# Row 17
$ws.Range("D17").Value = '2.054.68'
$ws.Range("E17").Value = '  -1.86%  '

# Row 18
$ws.Range("D18").Value = '37.754.51'
$ws.Range("E18").Value = '  +0.09%  '

# Row 19
$ws.Range("D19").Value = '6.18'
$ws.Range("E19").Value = '  -0.19%  '

# Row 20
$ws.Range("D20").Value = '71.09'
$ws.Range("E20").Value = '  +1.38%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0833'
$ws.Range("E21").Value = '  +1.37%  '

# Row 22
$ws.Range("D22").Value = '228.55'
$ws.Range("E22").Value = '  +0.76%  '

# Row 23
$ws.Range("E23").Value = '  +0.04%  '

# Row 24
$ws.Range("E24").Value = '  -1.24%  '

# Row 25
$ws.Range("E25").Value = '  -0.14%  '

# Row 26
$ws.Range("D26").Value = '170.88'
$ws.Range("E26").Value = '  +1.47%  '

# Row 27
$ws.Range("D27").Value = '0.139'
$ws.Range("E27").Value = '  +5.20%  '

# Row 28
$ws.Range("E28").Value = '  +1.11%  '

# Row 29
$ws.Range("E29").Value = '  +0.18%  '

# Row 30
$ws.Range("E30").Value = '  -1.98%  '

# Row 31
$ws.Range("E31").Value = '  +2.70%  '

# Row 32
$ws.Range("D32").Value = '4.70'
$ws.Range("E32").Value = '  +1.50%  '

# Row 33
$ws.Range("D33").Value = '0.0631'
$ws.Range("E33").Value = '  +1.47%  '

# Row 34
$ws.Range("D34").Value = '4.67'
$ws.Range("E34").Value = '  +2.31%  '

# Row 35
$ws.Range("E35").Value = '  -3.53%  '

# Row 36
$ws.Range("E36").Value = '  +0.40%  '

# Row 37
$ws.Range("D37").Value = '3.40'
$ws.Range("E37").Value = '  -1.25%  '

# Row 38
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  -0.09%  '

# Row 39
$ws.Range("E39").Value = '  -2.10%  '

# Row 40
$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D40").Value = '99.94'
$ws.Range("E40").Value = '  +3.70%  '

# Row 41
$ws.Range("B41").Value = 'Cronos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D41").Value = '0.0971'
$ws.Range("E41").Value = '  -1.95%  '

# Row 42
$ws.Range("E42").Value = '  -2.05%  '

# Row 43
$ws.Range("E43").Value = '  +1.01%  '

# Row 44
$ws.Range("D44").Value = '1.453.00'
$ws.Range("E44").Value = '  -1.12%  '

# Row 45
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = '16.44'
$ws.Range("E45").Value = '  +7.31%  '

# Row 46
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").Value = '1.16'
$ws.Range("E46").Value = '  -1.39%  '

# Row 47
$ws.Range("D47").Value = '4.20'
$ws.Range("E47").Value = '  +1.55%  '

# Row 48
$ws.Range("E48").Value = '  +1.21%  '

# Row 49
$ws.Range("D49").Value = '7.42'
$ws.Range("E49").Value = '  +1.64%  '

# Row 50
$ws.Range("E50").Value = '  -0.82%  '

# Row 51
$ws.Range("D51").Value = '2.269.93'
$ws.Range("E51").Value = '  -1.20%  '
